# Update Name of Algo
# This updates several recomputed numeric values produced by the
# RandomForest imputation algorithm on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = 6.485899999999999
$ws.Range("A3").Value  = -21.49600000000003
$ws.Range("B5").Value  = 4.666300000000003
$ws.Range("C5").Value  = -14.13380000000001
$ws.Range("E7").Value  = 11.6149
$ws.Range("C9").Value  = -11.71870000000001
$ws.Range("C11").Value = -13.64590000000001
$ws.Range("E11").Value = 13.4531
$ws.Range("A14").Value = -20.44009999999998
$ws.Range("A16").Value = -21.49740000000002
$ws.Range("B16").Value = 5.350300000000003
$ws.Range("C17").Value = -11.5311
$ws.Range("E19").Value = 13.8003
$ws.Range("A21").Value = -21.34030000000001
$ws.Range("C21").Value = -11.2464
$ws.Range("E21").Value = 12.44500000000001
$ws.Range("A23").Value = -21.57230000000002
$ws.Range("A25").Value = -22.59170000000003
